$d = $word.ActiveDocument

# Locate the paragraph that currently holds the "Framework keuze ..." text
# (it's the 8th paragraph in the document at the time of editing).
$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Framework keuze toegelicht*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq $null) {
    throw "Could not locate target paragraph"
}

$startPara = $d.Paragraphs.Item($targetIndex)
$lastPara  = $d.Paragraphs.Item($d.Paragraphs.Count)

# Remove every paragraph after the target one (they all get folded into the
# single rewritten paragraph), but keep the target paragraph's own mark so
# its pPr / identity survive.
if ($lastPara.Range.End -gt $startPara.Range.End) {
    $deleteRange = $d.Range($d.Paragraphs.Item($targetIndex + 1).Range.Start, $lastPara.Range.End)
    $deleteRange.Delete()
}

# Re-fetch the (now last) paragraph and replace its run content — leaving
# the paragraph mark (and therefore its w:p/w:pPr) untouched — with the new
# set of runs that make up the rewritten text.
$para = $d.Paragraphs.Item($d.Paragraphs.Count)
$contentRange = $d.Range($para.Range.Start, $para.Range.End - 1)

$newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>I</w:t></w:r><w:r><w:t>n dit wekelijkse coachgesprek hebben wij het voornamelijk gehad over de definitieve keuze van ons framework. In eerste instantie hadden wij namelijk besloten om te gaan werken met Zend, hier zijn wij echter op teruggekomen na tegenvallende ervaringen. Hierna zijn we een dag bezig geweest met Yii, ook dit was niet het framework wat wij zochten. Uiteindelijk hebben wij voor CodeIgniter gekozen en uitgelegd waarom.  De coach gaf aan dat we deze verandering wel moeten doorvoeren in onze rapporten (aanpassen en als een nieuwe versie opslaan)</w:t></w:r><w:r><w:t xml:space="preserve">. Daarna hebben we kort </w:t></w:r><w:r><w:t>u</w:t></w:r><w:r><w:t xml:space="preserve">itgelegd </w:t></w:r><w:r><w:t>waar de les ‘Rapporteren voor Techniek’</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>over gaat. Verder hebben we het gehad over het LDAP systeem, we hebben de coach verteld dat wij m.b.v. een LDAP-browser toegang hebben gekregen tot het LDAP systeem en hebben onze coach onze bevindingen medegedeeld. We hebben hem verteld dat LDAP vervuild is, het toetsen van gegevens is op dit moment zeer moeilijk. Dit komt voornamelijk doordat de rollen van docenten op verschillende manieren wordt aangegeven. Tenslotte hebben we afgesproken dat deze week (week 20) een verdiepingsweek is (in CodeIgniter, ons framework) en dat wij indien mogelijk van start gaan met het bouwen van de website.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$contentRange.InsertXML($newXml)

Write-Output "Rewrote paragraph $targetIndex; paragraphs now: $($d.Paragraphs.Count)"
